# Update header labels and data labels from Chinese placeholder strings
# ("字符串A0" ... "字符串标题A", etc.) to plain English equivalents
# ("StringA0" ... "StringA", etc.), then move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: column headers
$ws.Range("A1").Value = "StringA"
$ws.Range("B1").Value = "StringB"
$ws.Range("C1").Value = "StringC"

# Rows 2-11: data values for columns A, B, C
$letters = @("A", "B", "C")
for ($col = 1; $col -le 3; $col++) {
    $letter = $letters[$col - 1]
    for ($i = 0; $i -le 9; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, $col).Value = "String" + $letter + $i
    }
}

# Update the active selection to C2:C11
$ws.Range("C2:C11").Select()
